$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11
$ws.Range("A11").Value = 10154.81
$ws.Range("B11").Value = 10068.219999999999
$ws.Range("C11").Value = 305.24
$ws.Range("D11").Value = 307.86
$ws.Range("E11").Value = $false
$ws.Range("F11").Value = 0.86
$ws.Range("G11").Value = 42613.765474537038
$ws.Range("H11").Value = $true

# Row 12
$ws.Range("A12").Value = 10086.77
$ws.Range("B12").Value = 10154.81
$ws.Range("C12").Value = 307.68
$ws.Range("D12").Value = 305.63
$ws.Range("E12").Value = $false
$ws.Range("F12").Value = -0.67
$ws.Range("G12").Value = 42614.672696759262
$ws.Range("H12").Value = $false

# Row 13
$ws.Range("A13").Value = 10023.219999999999
$ws.Range("B13").Value = 10086.77
$ws.Range("C13").Value = 307.95999999999998
$ws.Range("D13").Value = 306.02
$ws.Range("E13").Value = $false
$ws.Range("F13").Value = -0.63
$ws.Range("G13").Value = 42615.750069444446
$ws.Range("H13").Value = $false
